# Insert one new daily price record for "Poroto verde" (Femacal de La Calera)
# at row 369, pushing the existing rows 369:456 down to 370:457.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 369 (shifts old row 369 -> 370, ..., old 456 -> 457).
$ws.Rows.Item(369).Insert()

$newRow = 369
$ws.Cells.Item($newRow, 1).Value2 = 3
$ws.Cells.Item($newRow, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value2 = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value2 = 44855
$ws.Cells.Item($newRow, 5).Value2 = 5
$ws.Cells.Item($newRow, 6).Value2 = 100112031
$ws.Cells.Item($newRow, 7).Value2 = "Poroto verde"
$ws.Cells.Item($newRow, 8).Value2 = "Magnum"
$ws.Cells.Item($newRow, 9).Value2 = "Primera"
$ws.Cells.Item($newRow, 10).Value2 = 85
$ws.Cells.Item($newRow, 11).Value2 = 31000
$ws.Cells.Item($newRow, 12).Value2 = 32000
$ws.Cells.Item($newRow, 13).Value2 = 31471
$ws.Cells.Item($newRow, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item($newRow, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value2 = 1259
$ws.Cells.Item($newRow, 17).Value2 = 25
$ws.Cells.Item($newRow, 18).Value2 = "Hortaliza"
